# epexspot_prices.xlsx update
#  - "Prix Spot": insert a new date column (09-nov) before column DN
#    (i.e. before the existing "01-oct." column), filled with "-" for the
#    hourly data rows (no data yet for that day).
#  - "Gaz" / "CO2": append a new row for 2025-11-07.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Prix Spot": insert column DN (118th column) -> everything from the
# old DN onward (01-oct. ... 31-oct.) shifts right by one column, and the
# freshly inserted column becomes the new "09-nov" header with "-" placeholders.
# ---------------------------------------------------------------------------
$wsPrix = $wb.Worksheets.Item("Prix Spot")

$newColIndex = 118   # column DN

$wsPrix.Columns.Item($newColIndex).Insert()

$wsPrix.Cells.Item(1, $newColIndex).Value = "09-nov"

for ($r = 2; $r -le 25; $r++) {
    $wsPrix.Cells.Item($r, $newColIndex).Value = "-"
}

# ---------------------------------------------------------------------------
# Sheet "Gaz": append row 146 for 2025-11-07.
# ---------------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")

$wsGaz.Range("A146").NumberFormat = "@"
$wsGaz.Range("A146").Value = "2025-11-07"
$wsGaz.Range("A146").Style = "Normal"
$wsGaz.Range("B146").Value = 29.74

# ---------------------------------------------------------------------------
# Sheet "CO2": append row 146 for 2025-11-07.
# ---------------------------------------------------------------------------
$wsCO2 = $wb.Worksheets.Item("CO2")

$wsCO2.Range("A146").NumberFormat = "@"
$wsCO2.Range("A146").Value = "2025-11-07"
$wsCO2.Range("A146").Style = "Normal"
$wsCO2.Range("B146").Value = 79.36
